# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Update the DAMSLTag (column I) and DialogAct
# (column J) values for the rows whose tags changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = @("ba", "Appreciation")
    8   = @("%", "Uninterpretable")
    10  = @("sv", "Statement-opinion")
    16  = @("sv", "Statement-opinion")
    18  = @("%", "Uninterpretable")
    24  = @("sd", "Statement-non-opinion")
    25  = @("sd", "Statement-non-opinion")
    26  = @("sv", "Statement-opinion")
    37  = @("ba", "Appreciation")
    41  = @("sv", "Statement-opinion")
    59  = @("sv", "Statement-opinion")
    61  = @("sd", "Statement-non-opinion")
    63  = @("ba", "Appreciation")
    64  = @("sv", "Statement-opinion")
    73  = @("ba", "Appreciation")
    90  = @("ba", "Appreciation")
    95  = @("aa", "Agree/Accept")
    105 = @("aa", "Agree/Accept")
    116 = @("aa", "Agree/Accept")
    138 = @("sv", "Statement-opinion")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
